$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated figures from the 2021-09-27 (EDT) data download.
# Columns B..M for rows 3..12 (age groups 12-15 through 90+).
$data = @(
    @(0,   383,  3727, 0,     130.3, 1361.5, 0,  0,  0,   0,   0,    0),
    @(76,  1536, 1509, 58.9,  519.8, 1086.2, 0,  0,  2,   0,   0,    1.4),
    @(237, 4091, 2783, 64.1,  629.2, 1105.3, 0,  3,  14,  0,   0.5,  5.6),
    @(507, 4499, 2797, 118.7, 804.5, 1543.6, 0,  0,  35,  0,   0,    19.3),
    @(511, 2729, 1952, 96.5,  661.7, 1401.3, 2,  11, 83,  0.4, 2.7,  59.6),
    @(309, 1189, 1047, 61.4,  461.6, 1098.9, 4,  20, 105, 0.8, 7.8,  110.2),
    @(281, 708,  669,  53.2,  494.6, 1023.2, 11, 34, 104, 2.1, 23.8, 159.1),
    @(211, 330,  370,  52.4,  450.1, 822.6,  16, 27, 80,  4,   36.8, 177.8),
    @(136, 206,  212,  77.6,  543.3, 976.5,  12, 20, 38,  6.8, 52.7, 175),
    @(61,  58,   78,   171.4, 449.6, 931.1,  5,  5,  18,  14,  38.8, 214.9)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 3 + $i
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $col = 2 + $j  # column B = 2
        $ws.Cells.Item($row, $col).Value = $vals[$j]
    }
}
